$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the filename referenced in column C (C1:C28) to the new .mat file
$ws.Range("C1:C28").Value = "painHealthyCoords_20190426.mat"

# Move the active selection to C34 (matches the saved cursor position)
$ws.Range("C34").Select()

# Best-effort: restore the window horizontal screen position (xWindow) to 0
try { $excel.Left = 0 } catch { }
try { $wb.Windows.Item(1).Left = 0 } catch { }
